$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsBdctba = $wb.Worksheets.Item("BDCTBA")

# Update the "About" sheet text describing the default border-adjustment behavior:
# the EPS now assumes carbon tax border adjustments do NOT apply across industries by default.
$wsAbout.Range("A9").Value = "By default, the EPS assumes carbon tax border adjustments do not apply across"
$wsAbout.Range("A10").Value = "industries."

# Flip the boolean lever value on the BDCTBA sheet to disable the carbon tax border adjustment.
$wsBdctba.Range("B2").Value = 1

# Restore cell selections on each sheet, leaving "BDCTBA" selected first and then
# re-activating "About" so it remains the active/visible tab, matching the saved workbook state.
$wsBdctba.Activate()
$wsBdctba.Range("B3").Select()

$wsAbout.Activate()
$wsAbout.Range("A10").Select()

$wb.Save()
